$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.063760280609131
$ws.Range("B1").Value = 2.37343955039978
$ws.Range("C1").Value = 6.499833106994629
$ws.Range("D1").Value = 2.248624801635742
$ws.Range("E1").Value = 1.293314337730408
